# Apply updated crypto price/volume data per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function SetTextValue {
    param($Range, [string]$Text)
    # Force the cell to store a literal text value (matches the source
    # inlineStr cells) instead of letting Excel auto-convert numeric-
    # looking strings into real numbers. Restore the original style
    # afterwards so no stray number-format style is left behind.
    $origStyle = $Range.Style
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = $origStyle
}

SetTextValue $ws.Range("D2") '29.909.58'
SetTextValue $ws.Range("E2") '  +0.68%  '
SetTextValue $ws.Range("D3") '1.634.08'
SetTextValue $ws.Range("E3") '  +0.99%  '
SetTextValue $ws.Range("D4") '0.999'
SetTextValue $ws.Range("E4") '  +0.75%  '
SetTextValue $ws.Range("D5") '215.34'
SetTextValue $ws.Range("E5") '  +1.20%  '
SetTextValue $ws.Range("E6") '  -0.09%  '
SetTextValue $ws.Range("E7") '  +0.84%  '
SetTextValue $ws.Range("D8") '28.72'
SetTextValue $ws.Range("E8") '  -1.77%  '
SetTextValue $ws.Range("D9") '0.260'
SetTextValue $ws.Range("E9") '  +0.58%  '
SetTextValue $ws.Range("E11") '  -0.97%  '
SetTextValue $ws.Range("D12") '1.868.04'
SetTextValue $ws.Range("E12") '  +1.11%  '
SetTextValue $ws.Range("D13") '1.633.96'
SetTextValue $ws.Range("E13") '  +1.28%  '
SetTextValue $ws.Range("D14") '0.586'
SetTextValue $ws.Range("E14") '  +2.94%  '
SetTextValue $ws.Range("D15") '9.51'
SetTextValue $ws.Range("E15") '  +5.17%  '
SetTextValue $ws.Range("D16") '29.910.15'
SetTextValue $ws.Range("E16") '  +0.69%  '
SetTextValue $ws.Range("D17") '3.85'
SetTextValue $ws.Range("E17") '  -1.68%  '
SetTextValue $ws.Range("D18") '65.58'
SetTextValue $ws.Range("E18") '  +2.16%  '
SetTextValue $ws.Range("D19") '241.32'
SetTextValue $ws.Range("E19") '  -0.29%  '
SetTextValue $ws.Range("E20") '  -1.07%  '
SetTextValue $ws.Range("D21") '0.999'
SetTextValue $ws.Range("E21") '  +0.60%  '
SetTextValue $ws.Range("E22") '  +1.61%  '
SetTextValue $ws.Range("D23") '4.15'
SetTextValue $ws.Range("E23") '  +0.86%  '
SetTextValue $ws.Range("E24") '  +2.64%  '
SetTextValue $ws.Range("D25") '157.79'
SetTextValue $ws.Range("E25") '  +0.90%  '
SetTextValue $ws.Range("D26") '15.53'
SetTextValue $ws.Range("E26") '  -1.00%  '
SetTextValue $ws.Range("E27") '  -1.13%  '
SetTextValue $ws.Range("D28") '6.64'
SetTextValue $ws.Range("E28") '  +0.59%  '
SetTextValue $ws.Range("E29") '  +0.56%  '
SetTextValue $ws.Range("E30") '  +0.16%  '
SetTextValue $ws.Range("E31") '  +1.66%  '
SetTextValue $ws.Range("E32") '  +2.01%  '
SetTextValue $ws.Range("D33") '3.19'
SetTextValue $ws.Range("E33") '  -0.73%  '
SetTextValue $ws.Range("D34") '1.425.12'
SetTextValue $ws.Range("E34") '  -0.05%  '
SetTextValue $ws.Range("E35") '  +3.60%  '
SetTextValue $ws.Range("E36") '  -2.48%  '
SetTextValue $ws.Range("E37") '  -3.29%  '
SetTextValue $ws.Range("E38") '  +0.12%  '
SetTextValue $ws.Range("E39") '  +0.28%  '
SetTextValue $ws.Range("D40") '75.64'
SetTextValue $ws.Range("E40") '  +8.06%  '
SetTextValue $ws.Range("D41") '0.558'
SetTextValue $ws.Range("E41") '  -0.05%  '
SetTextValue $ws.Range("E42") '  +0.64%  '
SetTextValue $ws.Range("E43") '  +1.34%  '
SetTextValue $ws.Range("D44") '0.0502'
SetTextValue $ws.Range("E44") '  -0.61%  '
SetTextValue $ws.Range("E45") '  +0.80%  '
SetTextValue $ws.Range("E46") '  +0.57%  '
SetTextValue $ws.Range("D47") '5.35'
SetTextValue $ws.Range("E47") '  -1.74%  '
SetTextValue $ws.Range("D48") '1.774.85'
SetTextValue $ws.Range("E48") '  +1.03%  '
SetTextValue $ws.Range("D49") '48.49'
SetTextValue $ws.Range("E49") '  -9.74%  '
SetTextValue $ws.Range("D50") '92.38'
SetTextValue $ws.Range("E50") '  +5.04%  '
SetTextValue $ws.Range("E51") '  +3.76%  '
